$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (K2:T2)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.083576666666666
$ws.Range("N2").Value = 9.250729999999999
$ws.Range("O2").Value = 0.2272509363535097
$ws.Range("P2").Value = 0.2272509363535097
$ws.Range("Q2").Value = 0.2197582861622222
$ws.Range("R2").Value = 1.97782457546
$ws.Range("S2").Value = 0.2272509363535097
$ws.Range("T2").Value = 0.2272509363535097

# Row 3 updates (M3:T3), K3/L3/N3 unchanged
$ws.Range("M3").Value = 6.453984666666667
$ws.Range("O3").Value = 0.4756405360586227
$ws.Range("P3").Value = 0.4756405360586227
$ws.Range("Q3").Value = 0.4599582765675556
$ws.Range("R3").Value = 4.139624489108
$ws.Range("S3").Value = 0.4756405360586227
$ws.Range("T3").Value = 0.4756405360586227

# Row 4 updates (M4:T4), K4/L4 unchanged
$ws.Range("M4").Value = 4.031477000000001
$ws.Range("N4").Value = 12.094431
$ws.Range("O4").Value = 0.2971085275878677
$ws.Range("P4").Value = 0.2971085275878677
$ws.Range("Q4").Value = 0.2873126151846667
$ws.Range("R4").Value = 2.585813536662001
$ws.Range("S4").Value = 0.2971085275878677
$ws.Range("T4").Value = 0.2971085275878677
